$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("111:111").Insert()
$ws.Range("A112:Q112").Copy($ws.Range("A111:Q111"))
$r = $ws.Range("Q111")
$r.Orientation = 0
$r.ReadingOrder = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.NumberFormat = "@"
Write-Host "done"
